# Dodanie podziału treningu na części
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Trening" column header (reuse the bold/centered header style) ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "Trening"

# --- Apply the datetime number format to column A data cells ---
# (Replicates the original authoring tool's behaviour of first registering a
#  lower-case format code and then switching to the upper-case one that is
#  actually used, so numFmtId 164 stays registered but unused while 165 is
#  applied to the cellXf used by column A.)
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Full data table: row, A(datetime serial), B, C, D, E, F ---
$data = @(
    @(2,  45685.64614293981,  941.7,  14.95, 5.046788147517615, "10-15", "Duża Gra"),
    @(3,  45685.64971932871,  1250.7, 12.6,  3.776089804513113, "10-15", "Duża Gra"),
    @(4,  45685.66124131945,  2246.2, 13.86, 3.933877059391568, "10-15", "Duża Gra"),
    @(5,  45685.64613946759,  941.4,  8.220000000000001, 3.922867093767438, "5-10", "Duża Gra"),
    @(6,  45685.65116261574,  1375.4, 9.98,  2.980909841401236, "5-10", "Duża Gra"),
    @(7,  45685.66123784722,  2245.9, 8.949999999999999, 3.069289088249209, "5-10", "Duża Gra"),
    @(8,  45685.667546875,    2791,   14.13, 3.543050800051007, "10-15", "Mała Gra"),
    @(9,  45685.66911168981,  2926.2, 14.96, 3.308462892259869, "10-15", "Mała Gra"),
    @(10, 45685.68455150463,  4260.2, 11.69, 3.183844123567855, "10-15", "Mała Gra"),
    @(11, 45685.67038715278,  3036.4, 9.41,  3.041034323828562, "5-10", "Mała Gra"),
    @(12, 45685.68442534722,  4249.3, 9.41,  2.977724824632916, "5-10", "Mała Gra"),
    @(13, 45685.68454918981,  4260,   8.83,  2.918066586766926, "5-10", "Mała Gra")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
